$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row values were reshuffled between existing records (dates/volumes/prices/origin swapped across rows).
# Set each affected cell to its new target value directly.

$ws.Range("D2").Value = 44362
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 6500
$ws.Range("L2").Value = 6500
$ws.Range("M2").Value = 6500
$ws.Range("N2").Value = '$/caja 36 atados'
$ws.Range("O2").Value = 'Región Metropolitana'
$ws.Range("P2").Value = 181
$ws.Range("Q2").Value = 36

$ws.Range("D3").Value = 44354
$ws.Range("J3").Value = 150
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 7000
$ws.Range("N3").Value = '$/caja 36 atados'
$ws.Range("O3").Value = 'Región del Maule'
$ws.Range("P3").Value = 194
$ws.Range("Q3").Value = 36

$ws.Range("D4").Value = 44355
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = 7000
$ws.Range("N4").Value = '$/caja 36 atados'
$ws.Range("O4").Value = 'Región Metropolitana'
$ws.Range("P4").Value = 194
$ws.Range("Q4").Value = 36

$ws.Range("D5").Value = 44372
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 7000
$ws.Range("N5").Value = '$/caja 36 atados'
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 194
$ws.Range("Q5").Value = 36

$ws.Range("D6").Value = 44342
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 7000
$ws.Range("N6").Value = '$/caja 36 atados'
$ws.Range("O6").Value = 'Región del Maule'
$ws.Range("P6").Value = 194
$ws.Range("Q6").Value = 36

$ws.Range("D7").Value = 44340
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 7000
$ws.Range("N7").Value = '$/caja 36 atados'
$ws.Range("O7").Value = 'Región del Maule'
$ws.Range("P7").Value = 194
$ws.Range("Q7").Value = 36

$ws.Range("D8").Value = 44386
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 6500
$ws.Range("L8").Value = 6500
$ws.Range("M8").Value = 6500
$ws.Range("N8").Value = '$/caja 36 atados'
$ws.Range("O8").Value = 'Región Metropolitana'
$ws.Range("P8").Value = 181
$ws.Range("Q8").Value = 36

$ws.Range("D9").Value = 44371
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 6500
$ws.Range("L9").Value = 6500
$ws.Range("M9").Value = 6500
$ws.Range("N9").Value = '$/caja 36 atados'
$ws.Range("O9").Value = 'Región Metropolitana'
$ws.Range("P9").Value = 181
$ws.Range("Q9").Value = 36

$ws.Range("D10").Value = 44348
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 7000
$ws.Range("N10").Value = '$/caja 36 atados'
$ws.Range("O10").Value = 'Región del Maule'
$ws.Range("P10").Value = 194
$ws.Range("Q10").Value = 36

$ws.Range("D11").Value = 44364
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("N11").Value = '$/caja 36 atados'
$ws.Range("O11").Value = 'Región Metropolitana'
$ws.Range("P11").Value = 194
$ws.Range("Q11").Value = 36

$ws.Range("D12").Value = 44358
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 7000
$ws.Range("N12").Value = '$/caja 36 atados'
$ws.Range("O12").Value = 'Región Metropolitana'
$ws.Range("P12").Value = 194
$ws.Range("Q12").Value = 36

$ws.Range("D13").Value = 44376
$ws.Range("J13").Value = 150
$ws.Range("K13").Value = 6500
$ws.Range("L13").Value = 6500
$ws.Range("M13").Value = 6500
$ws.Range("N13").Value = '$/caja 36 atados'
$ws.Range("O13").Value = 'Región Metropolitana'
$ws.Range("P13").Value = 181
$ws.Range("Q13").Value = 36

$ws.Range("D14").Value = 44357
$ws.Range("J14").Value = 150
$ws.Range("K14").Value = 6500
$ws.Range("L14").Value = 6500
$ws.Range("M14").Value = 6500
$ws.Range("N14").Value = '$/caja 20 docenas'
$ws.Range("O14").Value = 'Región del Maule'
$ws.Range("P14").Value = 6500
$ws.Range("Q14").Value = 1

$ws.Range("D15").Value = 44369
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("N15").Value = '$/caja 20 docenas'
$ws.Range("O15").Value = 'Región Metropolitana'
$ws.Range("P15").Value = 7000
$ws.Range("Q15").Value = 1

